$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Remove the obsolete "LPrgm2" test program row (row 27). Its shared string
# stays in the workbook's string table (unreferenced) - only the row goes.
$ws.Rows.Item(27).Delete()

# Insert a new row at 13 for "Artificial Intelligence" (alphabetically between
# "Annes Program 2" and "biology 3"), shifting rows below it down by one.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "Artificial Intelligence"

# Insert a new row at 28 for "LProgram1" (alphabetically between "KPS Genesis"
# and "Machine Learning" - that's where "Machine Learning" now sits after the
# delete+insert above).
$ws.Rows.Item(28).Insert()
$ws.Range("A28").Value = "LProgram1"
